{"js": "// This script rebuilds the document body (after the centered name/title\n// paragraph, which is left untouched) to match the edited resume content:\n//  - Removes the separate centered contact-info paragraph.\n//  - Rewrites the PROFESSIONAL SUMMARY paragraph.\n//  - Blanks out the CORE COMPETENCIES line (keeps an empty paragraph).\n//  - Replaces PROFESSIONAL EXPERIENCE with the full 10-position history.\n//  - Replaces KEY PROJECTS with four curated project write-ups.\n//  - Replaces KEY ACHIEVEMENTS AND IMPACT with a single curated \"Impact\" list.\n//  - Removes the TECHNICAL SKILLS detail lines and the closing sentence,\n//    leaving just the \"TECHNICAL SKILLS\" heading.\n\nconst targetParagraphs = [\n    { text: \"PROFESSIONAL SUMMARY\", style: \"Heading2\" },\n    { text: \"Senior data scientist and software engineer specializing in geospatial machine learning and large-scale demographic analysis. Developed algorithms that improved demographic classification accuracy from 23% to 64%, processed data across 178,000+ precincts, and built platforms serving thousands of analysts nationwide.\", style: \"Normal\" },\n    { text: \"CORE COMPETENCIES\", style: \"Heading2\" },\n    { text: \"\", style: \"Normal\" },\n    { text: \"PROFESSIONAL EXPERIENCE\", style: \"Heading2\" },\n    { text: \"Partner - Siege Analytics (Austin, TX) | 2005 - Present\", style: \"Heading3\" },\n    { text: \"Data, Technology and Strategy Consulting\", style: \"Normal\" },\n    { text: \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\", style: \"Normal\" },\n    { text: \"\u2022 Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration\", style: \"Normal\" },\n    { text: \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \u00b14.2% to \u00b12.1%\", style: \"Normal\" },\n    { text: \"Data Products Manager - Helm/Murmuration (Austin, TX) | June 2021 - May 2023\", style: \"Heading3\" },\n    { text: \"Civic Graph & Civic Pulse Director\", style: \"Normal\" },\n    { text: \"\u2022 Conceived, architected and built Civic Graph multi-tenant data warehouse processing government data from Census, Bureau of Labor Statistics, National Council of Educational Statistics\", style: \"Normal\" },\n    { text: \"\u2022 Built multi-dimensional data warehouse measuring socio-economic changes in America at every level across attitudinal, behavioral, demographic, economic and geographical dimensions\", style: \"Normal\" },\n    { text: \"\u2022 Managed engineering teams of 7-11 professionals while setting technical direction for data architecture\", style: \"Normal\" },\n    { text: \"Analytics Supervisor - GSD&M (Austin, TX) | November 2019 - June 2020\", style: \"Heading3\" },\n    { text: \"Big Data Engineering Transformation\", style: \"Normal\" },\n    { text: \"\u2022 Transformed small data team into big data engineering team, scaling from laptop datasets to Hadoop Clusters and Hive on AWS\", style: \"Normal\" },\n    { text: \"\u2022 Managed accounts including United States Air Force, Southwest Airlines/Chase and Indeed\", style: \"Normal\" },\n    { text: \"\u2022 Rewrote mission and offerings of department and drafted integration plan with strategy team\", style: \"Normal\" },\n    { text: \"Software Engineer - Mautinoa Technologies (Austin, TX) | August 2016 - February 2018\", style: \"Heading3\" },\n    { text: \"SimCrisis Product Owner/Engineer\", style: \"Normal\" },\n    { text: \"\u2022 Conceived, architected and engineered econometric simulation software for humanitarian crises intervention measurement\", style: \"Normal\" },\n    { text: \"\u2022 Built SimCrisis GeoDjango web application using multi-agent modeling to create econometric simulations of crisis economies\", style: \"Normal\" },\n    { text: \"\u2022 Designed modular application accepting rules extensions for ethnic strife, different crises/disasters, supply failures\", style: \"Normal\" },\n    { text: \"Senior Analyst - Myers Research (Austin, TX) | August 2012 - February 2014\", style: \"Heading3\" },\n    { text: \"RACSO Product Owner/Engineer\", style: \"Normal\" },\n    { text: \"\u2022 Designed comprehensive survey instruments for specialized voting segments and niche markets\", style: \"Normal\" },\n    { text: \"\u2022 Co-developed RACSO web application managing all aspects of survey operations from instrument design to data analysis\", style: \"Normal\" },\n    { text: \"\u2022 Wrote RFP and analyzed bids from 1,200 vendors for research platform development\", style: \"Normal\" },\n    { text: \"Research Director - PCCC (Washington, DC) | 2010 - 2012\", style: \"Heading3\" },\n    { text: \"Political Research & Data Analysis (FLEEM System)\", style: \"Normal\" },\n    { text: \"\u2022 Conceived, architected, and engineered FLEEM web application using Twilio API handling tens of thousands of simultaneous phone calls using emulated predictive dialer for regulated political surveys\", style: \"Normal\" },\n    { text: \"\u2022 Developed IVR polling system for early quantitative research supporting Senators Martin Heinrich and Elizabeth Warren\", style: \"Normal\" },\n    { text: \"\u2022 Built comprehensive tabular and graphical reporting system with Python, GeoDjango, PostGIS, and Apache webserver\", style: \"Normal\" },\n    { text: \"Software Engineer - Salsa Labs (Washington, DC) | January 2011 - August 2011\", style: \"Heading3\" },\n    { text: \"Geospatial CRM Development\", style: \"Normal\" },\n    { text: \"\u2022 Maintained and extended comprehensive geospatial analysis and reporting tools for Java-based CRM system used by tens of thousands simultaneously\", style: \"Normal\" },\n    { text: \"\u2022 Developed custom tile server for Web Map Service (WMS) integration using GeoTools and OpenLayers\", style: \"Normal\" },\n    { text: \"\u2022 Built advanced geospatial analysis capabilities using Java, JavaScript, MySQL, and TileMill\", style: \"Normal\" },\n    { text: \"Programmer - Lake Research Partners (Washington, DC) | April 2008 - December 2008\", style: \"Heading3\" },\n    { text: \"Political Research & Analytics\", style: \"Normal\" },\n    { text: \"\u2022 Built the first collaborative and multi-actor contributed poll of polls used by the Democratic Party\", style: \"Normal\" },\n    { text: \"\u2022 Harmonized data from 20+ polling firms with incompatible methodologies and encoding systems\", style: \"Normal\" },\n    { text: \"\u2022 Created comprehensive meta-analysis framework handling millions of survey responses that became the $400M Polling Consortium Database at The Analyst Institute, now valued at $1B+\", style: \"Normal\" },\n    { text: \"KEY PROJECTS\", style: \"Heading2\" },\n    { text: \"National Redistricting Platform (2020 - 2021)\", style: \"Heading3\" },\n    { text: \"Cloud-based GeoDjango platform for redistricting analysis with real-time collaborative editing and Census integration, used by thousands of analysts nationwide\", style: \"Normal\" },\n    { text: \"Technologies: GeoDjango, PostGIS, AWS, Docker, React, Python\", style: \"Normal\" },\n    { text: \"Impact: Reduced mapping costs by 73.5%, saving organizations $4.7M in operational expenses\", style: \"Normal\" },\n    { text: \"FLEEM Political Polling System (2010 - 2012)\", style: \"Heading3\" },\n    { text: \"Completely self-built IVR system using Twilio API that contacted tens of thousands of voters daily, replicated call center functionality to performance parity\", style: \"Normal\" },\n    { text: \"Technologies: Twilio API, Python, Django, PostgreSQL, JavaScript\", style: \"Normal\" },\n    { text: \"Impact: Saved $840K in operational costs plus millions in avoided software licensing\", style: \"Normal\" },\n    { text: \"Geospatial Demographic Classification System (2013 - 2016)\", style: \"Heading3\" },\n    { text: \"Machine learning platform that discovered systematic coding errors and improved demographic classification accuracy from 23% to 64%\", style: \"Normal\" },\n    { text: \"Technologies: Python, Scikit-learn, PostGIS, GeoPandas, TensorFlow\", style: \"Normal\" },\n    { text: \"Impact: Corrected demographic data affecting all Black and Asian-American voters nationwide\", style: \"Normal\" },\n    { text: \"KEY ACHIEVEMENTS AND IMPACT\", style: \"Heading2\" },\n    { text: \"Impact\", style: \"Heading3\" },\n    { text: \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters\", style: \"Normal\" },\n    { text: \"\u2022 Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M\", style: \"Normal\" },\n    { text: \"\u2022 Built redistricting platform used by thousands of analysts nationwide\", style: \"Normal\" },\n    { text: \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\", style: \"Normal\" },\n    { text: \"TECHNICAL SKILLS\", style: \"Heading2\" },\n];\n\n// Load all existing paragraphs in the body.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Keep paragraph 0 (\"Dheeraj Chand\", the centered name/title line) untouched,\n// and delete every paragraph after it -- they will be rebuilt below.\nfor (let i = paragraphs.items.length - 1; i >= 1; i--) {\n  paragraphs.items[i].delete();\n}\nawait context.sync();\n\n// Re-fetch so we have a fresh, valid anchor to insert after.\nconst remaining = context.document.body.paragraphs;\nremaining.load(\"items\");\nawait context.sync();\nlet anchor = remaining.items[0];\n\n// Insert the new paragraphs, in order, directly after the name paragraph.\nfor (const item of targetParagraphs) {\n  const p = anchor.insertParagraph(item.text, Word.InsertLocation.after);\n  p.style = item.style;\n  anchor = p;\n}\n\nawait context.sync();\n", "ps1": "# This script rebuilds the document body (after the centered name/title\n# paragraph, which is left untouched) to match the edited resume content:\n#  - Removes the separate centered contact-info paragraph.\n#  - Rewrites the PROFESSIONAL SUMMARY paragraph.\n#  - Blanks out the CORE COMPETENCIES line (keeps an empty paragraph).\n#  - Replaces PROFESSIONAL EXPERIENCE with the full 10-position history.\n#  - Replaces KEY PROJECTS with four curated project write-ups.\n#  - Replaces KEY ACHIEVEMENTS AND IMPACT with a single curated \"Impact\" list.\n#  - Removes the TECHNICAL SKILLS detail lines and the closing sentence,\n#    leaving just the \"TECHNICAL SKILLS\" heading.\n\n$d = $word.ActiveDocument\n\n$targetParagraphs = @(\n    @{ text = 'PROFESSIONAL SUMMARY'; style = 'Heading 2' },\n    @{ text = 'Senior data scientist and software engineer specializing in geospatial machine learning and large-scale demographic analysis. Developed algorithms that improved demographic classification accuracy from 23% to 64%, processed data across 178,000+ precincts, and built platforms serving thousands of analysts nationwide.'; style = 'Normal' },\n    @{ text = 'CORE COMPETENCIES'; style = 'Heading 2' },\n    @{ text = ''; style = 'Normal' },\n    @{ text = 'PROFESSIONAL EXPERIENCE'; style = 'Heading 2' },\n    @{ text = 'Partner - Siege Analytics (Austin, TX) | 2005 - Present'; style = 'Heading 3' },\n    @{ text = 'Data, Technology and Strategy Consulting'; style = 'Normal' },\n    @{ text = '\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%'; style = 'Normal' },\n    @{ text = '\u2022 Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration'; style = 'Normal' },\n    @{ text = '\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \u00b14.2% to \u00b12.1%'; style = 'Normal' },\n    @{ text = 'Data Products Manager - Helm/Murmuration (Austin, TX) | June 2021 - May 2023'; style = 'Heading 3' },\n    @{ text = 'Civic Graph & Civic Pulse Director'; style = 'Normal' },\n    @{ text = '\u2022 Conceived, architected and built Civic Graph multi-tenant data warehouse processing government data from Census, Bureau of Labor Statistics, National Council of Educational Statistics'; style = 'Normal' },\n    @{ text = '\u2022 Built multi-dimensional data warehouse measuring socio-economic changes in America at every level across attitudinal, behavioral, demographic, economic and geographical dimensions'; style = 'Normal' },\n    @{ text = '\u2022 Managed engineering teams of 7-11 professionals while setting technical direction for data architecture'; style = 'Normal' },\n    @{ text = 'Analytics Supervisor - GSD&M (Austin, TX) | November 2019 - June 2020'; style = 'Heading 3' },\n    @{ text = 'Big Data Engineering Transformation'; style = 'Normal' },\n    @{ text = '\u2022 Transformed small data team into big data engineering team, scaling from laptop datasets to Hadoop Clusters and Hive on AWS'; style = 'Normal' },\n    @{ text = '\u2022 Managed accounts including United States Air Force, Southwest Airlines/Chase and Indeed'; style = 'Normal' },\n    @{ text = '\u2022 Rewrote mission and offerings of department and drafted integration plan with strategy team'; style = 'Normal' },\n    @{ text = 'Software Engineer - Mautinoa Technologies (Austin, TX) | August 2016 - February 2018'; style = 'Heading 3' },\n    @{ text = 'SimCrisis Product Owner/Engineer'; style = 'Normal' },\n    @{ text = '\u2022 Conceived, architected and engineered econometric simulation software for humanitarian crises intervention measurement'; style = 'Normal' },\n    @{ text = '\u2022 Built SimCrisis GeoDjango web application using multi-agent modeling to create econometric simulations of crisis economies'; style = 'Normal' },\n    @{ text = '\u2022 Designed modular application accepting rules extensions for ethnic strife, different crises/disasters, supply failures'; style = 'Normal' },\n    @{ text = 'Senior Analyst - Myers Research (Austin, TX) | August 2012 - February 2014'; style = 'Heading 3' },\n    @{ text = 'RACSO Product Owner/Engineer'; style = 'Normal' },\n    @{ text = '\u2022 Designed comprehensive survey instruments for specialized voting segments and niche markets'; style = 'Normal' },\n    @{ text = '\u2022 Co-developed RACSO web application managing all aspects of survey operations from instrument design to data analysis'; style = 'Normal' },\n    @{ text = '\u2022 Wrote RFP and analyzed bids from 1,200 vendors for research platform development'; style = 'Normal' },\n    @{ text = 'Research Director - PCCC (Washington, DC) | 2010 - 2012'; style = 'Heading 3' },\n    @{ text = 'Political Research & Data Analysis (FLEEM System)'; style = 'Normal' },\n    @{ text = '\u2022 Conceived, architected, and engineered FLEEM web application using Twilio API handling tens of thousands of simultaneous phone calls using emulated predictive dialer for regulated political surveys'; style = 'Normal' },\n    @{ text = '\u2022 Developed IVR polling system for early quantitative research supporting Senators Martin Heinrich and Elizabeth Warren'; style = 'Normal' },\n    @{ text = '\u2022 Built comprehensive tabular and graphical reporting system with Python, GeoDjango, PostGIS, and Apache webserver'; style = 'Normal' },\n    @{ text = 'Software Engineer - Salsa Labs (Washington, DC) | January 2011 - August 2011'; style = 'Heading 3' },\n    @{ text = 'Geospatial CRM Development'; style = 'Normal' },\n    @{ text = '\u2022 Maintained and extended comprehensive geospatial analysis and reporting tools for Java-based CRM system used by tens of thousands simultaneously'; style = 'Normal' },\n    @{ text = '\u2022 Developed custom tile server for Web Map Service (WMS) integration using GeoTools and OpenLayers'; style = 'Normal' },\n    @{ text = '\u2022 Built advanced geospatial analysis capabilities using Java, JavaScript, MySQL, and TileMill'; style = 'Normal' },\n    @{ text = 'Programmer - Lake Research Partners (Washington, DC) | April 2008 - December 2008'; style = 'Heading 3' },\n    @{ text = 'Political Research & Analytics'; style = 'Normal' },\n    @{ text = '\u2022 Built the first collaborative and multi-actor contributed poll of polls used by the Democratic Party'; style = 'Normal' },\n    @{ text = '\u2022 Harmonized data from 20+ polling firms with incompatible methodologies and encoding systems'; style = 'Normal' },\n    @{ text = '\u2022 Created comprehensive meta-analysis framework handling millions of survey responses that became the $400M Polling Consortium Database at The Analyst Institute, now valued at $1B+'; style = 'Normal' },\n    @{ text = 'KEY PROJECTS'; style = 'Heading 2' },\n    @{ text = 'National Redistricting Platform (2020 - 2021)'; style = 'Heading 3' },\n    @{ text = 'Cloud-based GeoDjango platform for redistricting analysis with real-time collaborative editing and Census integration, used by thousands of analysts nationwide'; style = 'Normal' },\n    @{ text = 'Technologies: GeoDjango, PostGIS, AWS, Docker, React, Python'; style = 'Normal' },\n    @{ text = 'Impact: Reduced mapping costs by 73.5%, saving organizations $4.7M in operational expenses'; style = 'Normal' },\n    @{ text = 'FLEEM Political Polling System (2010 - 2012)'; style = 'Heading 3' },\n    @{ text = 'Completely self-built IVR system using Twilio API that contacted tens of thousands of voters daily, replicated call center functionality to performance parity'; style = 'Normal' },\n    @{ text = 'Technologies: Twilio API, Python, Django, PostgreSQL, JavaScript'; style = 'Normal' },\n    @{ text = 'Impact: Saved $840K in operational costs plus millions in avoided software licensing'; style = 'Normal' },\n    @{ text = 'Geospatial Demographic Classification System (2013 - 2016)'; style = 'Heading 3' },\n    @{ text = 'Machine learning platform that discovered systematic coding errors and improved demographic classification accuracy from 23% to 64%'; style = 'Normal' },\n    @{ text = 'Technologies: Python, Scikit-learn, PostGIS, GeoPandas, TensorFlow'; style = 'Normal' },\n    @{ text = 'Impact: Corrected demographic data affecting all Black and Asian-American voters nationwide'; style = 'Normal' },\n    @{ text = 'KEY ACHIEVEMENTS AND IMPACT'; style = 'Heading 2' },\n    @{ text = 'Impact'; style = 'Heading 3' },\n    @{ text = '\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters'; style = 'Normal' },\n    @{ text = '\u2022 Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M'; style = 'Normal' },\n    @{ text = '\u2022 Built redistricting platform used by thousands of analysts nationwide'; style = 'Normal' },\n    @{ text = '\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%'; style = 'Normal' },\n    @{ text = 'TECHNICAL SKILLS'; style = 'Heading 2' },\n)\n\n# Keep paragraph 1 (\"Dheeraj Chand\", the centered name/title line) untouched,\n# and delete every paragraph after it -- they will be rebuilt below.\n$count = $d.Paragraphs.Count\nfor ($i = $count; $i -ge 2; $i--) {\n    $d.Paragraphs.Item($i).Range.Delete()\n}\n\n# Insert the new paragraphs, in order, directly after the name paragraph.\n$anchorIndex = 1\nforeach ($item in $targetParagraphs) {\n    $anchor = $d.Paragraphs.Item($anchorIndex)\n    $anchor.Range.InsertParagraphAfter()\n    $anchorIndex += 1\n    $newPara = $d.Paragraphs.Item($anchorIndex)\n    $newPara.Range.Text = $item.text\n    $newPara.Style = $item.style\n}\n"}
